$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Break the existing E12:E18 merge so the individual cells in that column
#    can be written to independently.
# ---------------------------------------------------------------------------
$ws.Range("E12:E18").UnMerge()

# ---------------------------------------------------------------------------
# 2. Add the five new "hour" rows (22-26) in column A, copying the formatting
#    (bold 14pt font, vertical-center alignment -> style used by A2:A21) from
#    the row right above so the new header cells match the existing ones.
# ---------------------------------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A22:A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A22").Value = "hour 21"
$ws.Range("A23").Value = "hour 22"
$ws.Range("A24").Value = "hour 23"
$ws.Range("A25").Value = "hour 24"
$ws.Range("A26").Value = "hour 25"

# ---------------------------------------------------------------------------
# 3. Replace the old "week 3" notes cell (E12, previously "1. Create item")
#    with the long Angular task list, and give the new "week 3" notes entry
#    (E25) the "organize the project" text.
# ---------------------------------------------------------------------------
$angularNotes = "1. Create assignment 1 project using Angular" + [char]10 +
  "2. Install Visual Studio Code and useful extentions" + [char]10 +
  "3. Create component for a item card " + [char]10 +
  "4. Create component for item list " + [char]10 +
  "5. Add css to item card using bootstrap" + [char]10 +
  "6. create service to make HTTP calls to retrieve data for each item card" + [char]10 +
  "7. Modify HTTP data using Pipe" + [char]10 +
  "8. Create routes for items in the navigation menu" + [char]10 +
  "9. Create new component for displaying details of each item" + [char]10 +
  "10. create next page button using event binding" + [char]10 +
  "11. Apply routing features to display different list of items based on different page urls" + [char]10 +
  "12. Apply defult picture for items with no image property" + [char]10 +
  "13. Create a new component for adding new items" + [char]10 +
  "14. Create a form for the new component using template driven approach" + [char]10 +
  "15. Apply validations to the form" + [char]10 +
  "16. Create animations and error notifacations that will display when inputs is not valid" + [char]10 +
  "17. Create a new component for user registration" + [char]10 +
  "18. Create a form for the new component using reactive form approach" + [char]10 +
  "19. apply validations" + [char]10 +
  "20. Create animations and error notifacations that will display when inputs is not valid using alertify" + [char]10 +
  "21. Organize the code for the form using 'FormBuilder'" + [char]10 +
  "22. Create a new service for adding new account information from the registration form to the local storage" + [char]10

$organizeNotes = "1. Organize the project, removing identity features and other unnecessary code that were auto-generated by visual studio 2019" + [char]10 +
  "2. Fix bootstrap problem and other startup issues after the operation"

$ws.Range("E12").Value = $angularNotes
$ws.Range("E25").Value = $organizeNotes

# ---------------------------------------------------------------------------
# 4. Give the whole E12:E26 block (the merged "week 3" notes column plus the
#    new "hour 21..25" rows) the same left/top-aligned, wrapped formatting.
# ---------------------------------------------------------------------------
$notesRange = $ws.Range("E12:E26")
$notesRange.HorizontalAlignment = -4131
$notesRange.VerticalAlignment = -4160
$notesRange.WrapText = $true

# ---------------------------------------------------------------------------
# 5. Re-merge E12:E24 (week-3 notes, now spanning the new hour rows) and the
#    new E25:E26 pair.
# ---------------------------------------------------------------------------
$ws.Range("E12:E24").Merge()
$ws.Range("E25:E26").Merge()

# ---------------------------------------------------------------------------
# 6. Every data row is now a uniform 60pt tall (instead of the old assorted
#    heights).
# ---------------------------------------------------------------------------
$ws.Range("A1:A26").RowHeight = 60

# ---------------------------------------------------------------------------
# 7. Update the view: scroll so row 19 / column C is at the top-left and the
#    active selection is D25 (matching the saved sheetView in the workbook).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D25").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 3
